$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Re-order the data-set comparison table from
#   A=Data Set, B=Read Length, C=Total Reads, D=Error Rate*,
#   E=Reads with Adapters*, F=Adapter Bases*
# to
#   A=Data Set, B=Error Rate, C=Read Length, D=Total Reads,
#   E=Reads w/ Adapters, F=Adapter Bases
# and add a new "GM12878 WGBS" row worth of data in row 5 (columns B,E,F).
# ---------------------------------------------------------------------------

# --- capture the original values (rows 2-4) before anything is overwritten
$readLength = @(125, 125, 125)
$totalReads = @(781923, 780899, 782237)
$errorRate  = @(0.002, 0.0060000000000000001, 0.012)
$readsWithAdapters = @(325982, 325105, 325860)
$adapterBases      = @(12447262, 12416861, 12464235)

# --- header row ----------------------------------------------------------
$ws.Range("A1").Value = "Data Set"

$ws.Range("B1").Value = "Error Rate"
$ws.Range("C1").Value = "Read Length"

$ws.Range("D1").Value = "Total Reads"
# D1 needs the "no center" bordered header look (same as A1) instead of the
# centered bordered header look it has today.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = "Total Reads"

$ws.Range("E1").Value = "Reads w/ Adapters"
$ws.Range("F1").Value = "Adapter Bases"

# --- column B : Error Rate (numeric, percent, centered) ------------------
for ($i = 0; $i -lt 3; $i++) {
    $row = 2 + $i
    $c = $ws.Range("B$row")
    $c.Value = $errorRate[$i]
    $c.NumberFormat = "0.00%"
    $c.HorizontalAlignment = -4108
}
$ws.Range("B5").Value = "Unknown"
$ws.Range("B5").HorizontalAlignment = -4108

# --- column C : Read Length (numeric, general, centered) -----------------
for ($i = 0; $i -lt 3; $i++) {
    $row = 2 + $i
    $c = $ws.Range("C$row")
    $c.Value = $readLength[$i]
    $c.HorizontalAlignment = -4108
}
$ws.Range("C5").Value = 125
$ws.Range("C5").HorizontalAlignment = -4108

# --- column D : Total Reads (numeric, #,##0, bordered flag) --------------
for ($i = 0; $i -lt 3; $i++) {
    $row = 2 + $i
    $c = $ws.Range("D$row")
    $c.Value = $totalReads[$i]
    $c.NumberFormat = "#,##0"
    $c.Borders.Item(9).LineStyle = 1
    $c.Borders.Item(9).LineStyle = -4142
}
$ws.Range("D5").Value = 1000000
$ws.Range("D5").NumberFormat = "#,##0"
$ws.Range("D5").Borders.Item(9).LineStyle = 1
$ws.Range("D5").Borders.Item(9).LineStyle = -4142

# --- column E : Reads w/ Adapters (numeric, #,##0) ------------------------
for ($i = 0; $i -lt 3; $i++) {
    $row = 2 + $i
    $c = $ws.Range("E$row")
    $c.Value = $readsWithAdapters[$i]
    $c.NumberFormat = "#,##0"
}
$ws.Range("E5").Value = "16,999*"
$ws.Range("E5").NumberFormat = "#,##0"
$ws.Range("E5").HorizontalAlignment = -4152

# --- column F : Adapter Bases (numeric, #,##0) ----------------------------
for ($i = 0; $i -lt 3; $i++) {
    $row = 2 + $i
    $c = $ws.Range("F$row")
    $c.Value = $adapterBases[$i]
    $c.NumberFormat = "#,##0"
}
$ws.Range("F5").Value = "1,020,017*"
$ws.Range("F5").NumberFormat = "#,##0"
$ws.Range("F5").HorizontalAlignment = -4152

# --- column widths ---------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 9.998697916666666
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(5).ColumnWidth = 16.330729166666668

# --- selection --------------------------------------------------------------
$ws.Range("F6").Select() | Out-Null
